# Auto-generated edit script: updates market price / profit columns (H-N)
# across multiple sheets, per the scheduled market-data refresh diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 8052689
$ws.Range("J138").Value = 8068715.5
$ws.Range("L138").Value = 24206146.5
$ws.Range("N138").Value = -24216426.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 127129.125
$ws.Range("I2").Value = 145076.14
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 145076.14
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -144963.14
$ws.Range("N2").Value = -1726

$ws.Range("H32").Value = 30769.797
$ws.Range("I32").Value = 8916.888999999999
$ws.Range("J32").Value = 82526.69
$ws.Range("K32").Value = 8916.888999999999
$ws.Range("L32").Value = 82526.69
$ws.Range("M32").Value = -8629.888999999999
$ws.Range("N32").Value = -83100.69

$ws.Range("H116").Value = 127129.125
$ws.Range("I116").Value = 145076.14
$ws.Range("J116").Value = 1500
$ws.Range("K116").Value = 145076.14
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = -142782.14
$ws.Range("N116").Value = -6088

$ws.Range("H132").Value = 1734.6666
$ws.Range("I132").Value = 1622.0677
$ws.Range("J132").Value = 2399
$ws.Range("K132").Value = 4866.203100000001
$ws.Range("L132").Value = 7197
$ws.Range("M132").Value = -2336.203100000001
$ws.Range("N132").Value = -12257

$ws.Range("H139").Value = 59500
$ws.Range("J139").Value = 59500
$ws.Range("L139").Value = 59500
$ws.Range("N139").Value = -69780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 127129.125
$ws.Range("I3").Value = 145076.14
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 145076.14
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = -144962.14
$ws.Range("N3").Value = -1728

$ws.Range("H13").Value = 37500
$ws.Range("J13").Value = 37500
$ws.Range("L13").Value = 37500
$ws.Range("N13").Value = -37836

$ws.Range("H134").Value = 2805.0667
$ws.Range("I134").Value = 1879.3158
$ws.Range("J134").Value = 7830.5713
$ws.Range("K134").Value = 5637.9474
$ws.Range("L134").Value = 23491.7139
$ws.Range("M134").Value = -3102.9474
$ws.Range("N134").Value = -28561.7139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1891.4
$ws.Range("I16").Value = 1804.2
$ws.Range("J16").Value = 1978.6
$ws.Range("K16").Value = 1804.2
$ws.Range("L16").Value = 1978.6
$ws.Range("M16").Value = -1517.2
$ws.Range("N16").Value = -2552.6

$ws.Range("H31").Value = 5098.123
$ws.Range("I31").Value = 2175.9546
$ws.Range("J31").Value = 6934.914
$ws.Range("K31").Value = 2175.9546
$ws.Range("L31").Value = 6934.914
$ws.Range("M31").Value = -1880.9546
$ws.Range("N31").Value = -7524.914

$ws.Range("H34").Value = 5098.123
$ws.Range("I34").Value = 2175.9546
$ws.Range("J34").Value = 6934.914
$ws.Range("K34").Value = 2175.9546
$ws.Range("L34").Value = 6934.914
$ws.Range("M34").Value = -1973.9546
$ws.Range("N34").Value = -7338.914

$ws.Range("H113").Value = 1891.4
$ws.Range("I113").Value = 1804.2
$ws.Range("J113").Value = 1978.6
$ws.Range("K113").Value = 1804.2
$ws.Range("L113").Value = 1978.6
$ws.Range("M113").Value = 365.8
$ws.Range("N113").Value = -6318.6

$ws.Range("H134").Value = 42859932
$ws.Range("I134").Value = 83334410
$ws.Range("J134").Value = 21742810
$ws.Range("K134").Value = 250003230
$ws.Range("L134").Value = 65228430
$ws.Range("M134").Value = -250000695
$ws.Range("N134").Value = -65233500

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2021.5807
$ws.Range("I5").Value = 921.7646999999999
$ws.Range("J5").Value = 3357.0715
$ws.Range("K5").Value = 2765.2941
$ws.Range("L5").Value = 10071.2145
$ws.Range("M5").Value = -2653.2941
$ws.Range("N5").Value = -10295.2145

$ws.Range("H107").Value = 557.7273
$ws.Range("I107").Value = 576.61536
$ws.Range("J107").Value = 530.44446
$ws.Range("K107").Value = 1729.84608
$ws.Range("L107").Value = 1591.33338
$ws.Range("M107").Value = 190.15392
$ws.Range("N107").Value = -5431.33338

$ws.Range("H113").Value = 5348593.5
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 5348593.5
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 16045780.5
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -16050120.5

$ws.Range("H118").Value = 3000
$ws.Range("I118").Value = 3000
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 9000
$ws.Range("L118").Value = 0
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = -7757

$ws.Range("H125").Value = 2925.3845

$ws.Range("H134").Value = 6877.4707
$ws.Range("I134").Value = 3146.4443
$ws.Range("J134").Value = 11074.875
$ws.Range("K134").Value = 9439.332900000001
$ws.Range("L134").Value = 33224.625
$ws.Range("M134").Value = -4369.332900000001
$ws.Range("N134").Value = -43364.625

$ws.Range("H135").Value = 2021.5807
$ws.Range("I135").Value = 921.7646999999999
$ws.Range("J135").Value = 3357.0715
$ws.Range("K135").Value = 8295.882299999999
$ws.Range("L135").Value = 30213.6435
$ws.Range("M135").Value = -5760.882299999999
$ws.Range("N135").Value = -35283.6435

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1228.4667
$ws.Range("I113").Value = 878.25
$ws.Range("J113").Value = 2629.3333
$ws.Range("K113").Value = 878.25
$ws.Range("L113").Value = 2629.3333
$ws.Range("M113").Value = 1291.75
$ws.Range("N113").Value = -6969.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 683.94116
$ws.Range("I46").Value = 644.36365
$ws.Range("J46").Value = 756.5
$ws.Range("K46").Value = 644.36365
$ws.Range("L46").Value = 756.5
$ws.Range("M46").Value = -456.36365
$ws.Range("N46").Value = -1132.5

$ws.Range("H117").Value = 32000
$ws.Range("J117").Value = 32000
$ws.Range("L117").Value = 32000
$ws.Range("N117").Value = -41178

$ws.Range("H122").Value = 2520.5264
$ws.Range("I122").Value = 1516.8182
$ws.Range("J122").Value = 3900.625
$ws.Range("K122").Value = 4550.4546
$ws.Range("L122").Value = 11701.875
$ws.Range("M122").Value = -2100.4546
$ws.Range("N122").Value = -16601.875

$ws.Range("H136").Value = 7691.294
$ws.Range("I136").Value = 3049.6667
$ws.Range("J136").Value = 12913.125
$ws.Range("K136").Value = 9149.000100000001
$ws.Range("L136").Value = 38739.375
$ws.Range("M136").Value = -6599.000100000001
$ws.Range("N136").Value = -43839.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 588.7273
$ws.Range("I100").Value = 588.7273
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1177.4546
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -636.4546

$ws.Range("H132").Value = 2391.75
$ws.Range("I132").Value = 2152.5637
$ws.Range("J132").Value = 3853.4443
$ws.Range("K132").Value = 6457.6911
$ws.Range("L132").Value = 11560.3329
$ws.Range("M132").Value = -3927.6911
$ws.Range("N132").Value = -16620.3329
